$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new row of data: A3 = 4381
$ws.Range("A3").Value = 4381

# Update selection to match the new active cell
$ws.Range("A3").Select()
